# Update the want-to-go count values in column F across the four
# worksheets (Exhibitions, Performances, Local Life, All Types) to
# reflect the refreshed data snapshot referenced by the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item(1)
$wsPerformances = $wb.Worksheets.Item(2)
$wsLocalLife = $wb.Worksheets.Item(3)
$wsAllTypes = $wb.Worksheets.Item(4)

# Sheet 1
$wsExhibitions.Range("F3").Value  = 27022
$wsExhibitions.Range("F5").Value  = 645
$wsExhibitions.Range("F7").Value  = 566
$wsExhibitions.Range("F8").Value  = 234
$wsExhibitions.Range("F10").Value = 475
$wsExhibitions.Range("F16").Value = 65
$wsExhibitions.Range("F17").Value = 1617
$wsExhibitions.Range("F19").Value = 548
$wsExhibitions.Range("F21").Value = 456

# Sheet 2
$wsPerformances.Range("F2").Value  = 4523
$wsPerformances.Range("F8").Value  = 41
$wsPerformances.Range("F11").Value = 453
$wsPerformances.Range("F18").Value = 28

# Sheet 3
$wsLocalLife.Range("F2").Value = 5175
$wsLocalLife.Range("F3").Value = 271

# Sheet 4
$wsAllTypes.Range("F3").Value  = 5175
$wsAllTypes.Range("F4").Value  = 271
$wsAllTypes.Range("F5").Value  = 27022
$wsAllTypes.Range("F6").Value  = 4523
$wsAllTypes.Range("F9").Value  = 645
$wsAllTypes.Range("F15").Value = 41
$wsAllTypes.Range("F18").Value = 453
$wsAllTypes.Range("F19").Value = 566
$wsAllTypes.Range("F22").Value = 234
$wsAllTypes.Range("F24").Value = 475
$wsAllTypes.Range("F33").Value = 65
$wsAllTypes.Range("F35").Value = 1617
$wsAllTypes.Range("F37").Value = 548
$wsAllTypes.Range("F38").Value = 28
$wsAllTypes.Range("F40").Value = 456

$wb.Save()
